$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "diego"
$ws.Range("C2").Value = "fnjrdrn@gmail"
$ws.Range("G2").Value = 2
$ws.Range("H2").Value = 95.45
